$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 307, shifting existing rows 307-323 down to 308-324.
$ws.Rows(307).Insert()

# Populate the new row 307 with the new weekly price record.
$ws.Range("A307").Value = 5
$ws.Range("B307").Value = 'Macroferia Regional de Talca'
$ws.Range("C307").Value = 'Maule'
$ws.Range("D307").Value = 44753
$ws.Range("E307").Value = 7
$ws.Range("F307").Value = 100112003
$ws.Range("G307").Value = 'Ajo'
$ws.Range("H307").Value = 'Chino'
$ws.Range("I307").Value = 'Primera'
$ws.Range("J307").Value = 300
$ws.Range("K307").Value = 21000
$ws.Range("L307").Value = 21000
$ws.Range("M307").Value = 21000
$ws.Range("N307").Value = '$/malla 10 kilos'
$ws.Range("O307").Value = 'China'
$ws.Range("P307").Value = 2100
$ws.Range("Q307").Value = 10
$ws.Range("R307").Value = 'Hortaliza'
